# Edit script replicating the commit "Add files via upload"
# Changes:
#  1. Cover!C12: date serial 44084 -> literal text "13/10/202" (new shared string)
#  2. sa_all_tot: update raw F2:F49 values (dependent formulas G/I and the
#     sa_agg sheet's linked D column recalc automatically)
#  3. sa_all_tot: replace the plain cell selection with a frozen-pane view
#     (freeze row 1 + column A, i.e. top-left unfrozen cell B2)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Cover sheet: C12 becomes literal text instead of a date value
# ---------------------------------------------------------------------------
$coverWs = $wb.Worksheets.Item("Cover")
$coverWs.Range("C12").Value = "13/10/202"

# ---------------------------------------------------------------------------
# 2) sa_all_tot: refresh the raw stock-addition figures in column F
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("sa_all_tot")

$ws.Range("F2").Value = 101272075.7867339
$ws.Range("F3").Value = 109021675.1835863
$ws.Range("F4").Value = 105802274.0149852
$ws.Range("F5").Value = 37979511.53022062
$ws.Range("F6").Value = 550940375.47330797
$ws.Range("F7").Value = 393306803.89212358
$ws.Range("F8").Value = 67714485.879250303
$ws.Range("F9").Value = 13826419724.27739
$ws.Range("F10").Value = 11579109.315412
$ws.Range("F11").Value = 90337498.874392003
$ws.Range("F12").Value = 700117418.71123314
$ws.Range("F13").Value = 52398711.876658797
$ws.Range("F14").Value = 12108405.31857335
$ws.Range("F15").Value = 281205748.08100861
$ws.Range("F16").Value = 81850768.140337244
$ws.Range("F17").Value = 429102918.14303857
$ws.Range("F18").Value = 343077724.11224473
$ws.Range("F19").Value = 75384562.434790298
$ws.Range("F20").Value = 20150230.748340368
$ws.Range("F21").Value = 34232179.177534357
$ws.Range("F22").Value = 376647416.63428497
$ws.Range("F23").Value = 48438423.383421257
$ws.Range("F24").Value = 1635958539.3553841
$ws.Range("F25").Value = 386988961.37466788
$ws.Range("F26").Value = 619006497.90629971
$ws.Range("F27").Value = 362022452.31414533
$ws.Range("F28").Value = 14164988.95837347
$ws.Range("F29").Value = 9891764.5639806371
$ws.Range("F30").Value = 12883649.050343711
$ws.Range("F31").Value = 2047962.915295037
$ws.Range("F32").Value = 262492713.06024149
$ws.Range("F33").Value = 115124385.7170886
$ws.Range("F34").Value = 73018685.88128984
$ws.Range("F35").Value = 375166197.6098035
$ws.Range("F36").Value = 107619974.0045954
$ws.Range("F37").Value = 101049621.6914745
$ws.Range("F38").Value = 531933323.90010369
$ws.Range("F39").Value = 108692573.3013818
$ws.Range("F40").Value = 13924634.445961369
$ws.Range("F41").Value = 38335311.195605762
$ws.Range("F42").Value = 371948812.21830678
$ws.Range("F43").Value = 2134949211.4876561
$ws.Range("F44").Value = 1728577979.9677351
$ws.Range("F45").Value = 140108473.83452389
$ws.Range("F46").Value = 920271693.76920414
$ws.Range("F47").Value = 745233544.2525481
$ws.Range("F48").Value = 1615799926.6286709
$ws.Range("F49").Value = 119797068.6323203

# ---------------------------------------------------------------------------
# 3) sa_all_tot: freeze the header row + first column (pane split at B2)
# ---------------------------------------------------------------------------
[void]$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Restore the original active sheet (Cover) so the saved file keeps its
# original tab-selection state.
[void]$coverWs.Activate()
